# Auto-generated script to apply cell value updates per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 205.2
$ws.Range("I4").Value = 205.2
$ws.Range("K4").Value = 205.2
$ws.Range("M4").Value = -91.19999999999999
$ws.Range("H17").Value = 587.80884
$ws.Range("J17").Value = 589.8806
$ws.Range("L17").Value = 1769.6418
$ws.Range("N17").Value = -2105.6418
$ws.Range("H51").Value = 102507000
$ws.Range("I51").Value = 205005800
$ws.Range("J51").Value = 8199.799999999999
$ws.Range("K51").Value = 205005800
$ws.Range("L51").Value = 8199.799999999999
$ws.Range("M51").Value = -205005316
$ws.Range("N51").Value = -9167.799999999999
$ws.Range("H70").Value = 120223.336
$ws.Range("J70").Value = 11444.286
$ws.Range("L70").Value = 34332.858
$ws.Range("N70").Value = -34872.858
$ws.Range("H73").Value = 120223.336
$ws.Range("J73").Value = 11444.286
$ws.Range("L73").Value = 34332.858
$ws.Range("N73").Value = -36204.858
$ws.Range("H86").Value = 5071.1577
$ws.Range("I86").Value = 4606.5557
$ws.Range("J86").Value = 5489.3
$ws.Range("K86").Value = 4606.5557
$ws.Range("L86").Value = 5489.3
$ws.Range("M86").Value = -3483.5557
$ws.Range("N86").Value = -7735.3
$ws.Range("H89").Value = 5071.1577
$ws.Range("I89").Value = 4606.5557
$ws.Range("J89").Value = 5489.3
$ws.Range("K89").Value = 23032.7785
$ws.Range("L89").Value = 27446.5
$ws.Range("M89").Value = -17416.7785
$ws.Range("N89").Value = -38678.5
$ws.Range("H121").Value = 2900
$ws.Range("J121").Value = 2900
$ws.Range("L121").Value = 8700
$ws.Range("N121").Value = -12194
$ws.Range("H137").Value = 2256
$ws.Range("I137").Value = 2145.875
$ws.Range("K137").Value = 6437.625
$ws.Range("M137").Value = -3887.625
$ws.Range("H141").Value = 4619.1333
$ws.Range("I141").Value = 3649
$ws.Range("J141").Value = 6074.3335
$ws.Range("K141").Value = 10947
$ws.Range("L141").Value = 18223.0005
$ws.Range("M141").Value = -5767
$ws.Range("N141").Value = -28583.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2602.2222
$ws.Range("I2").Value = 2593.7646
$ws.Range("K2").Value = 2593.7646
$ws.Range("M2").Value = -2480.7646
$ws.Range("H45").Value = 4368.5415
$ws.Range("I45").Value = 2438
$ws.Range("K45").Value = 2438
$ws.Range("M45").Value = -2061
$ws.Range("H63").Value = 238816.67
$ws.Range("I63").Value = 6666.6665
$ws.Range("J63").Value = 470966.66
$ws.Range("K63").Value = 6666.6665
$ws.Range("L63").Value = 470966.66
$ws.Range("M63").Value = -5980.6665
$ws.Range("N63").Value = -472338.66
$ws.Range("H66").Value = 238816.67
$ws.Range("I66").Value = 6666.6665
$ws.Range("J66").Value = 470966.66
$ws.Range("K66").Value = 33333.3325
$ws.Range("L66").Value = 2354833.3
$ws.Range("M66").Value = -29901.3325
$ws.Range("N66").Value = -2361697.3
$ws.Range("H74").Value = 925.9655
$ws.Range("I74").Value = 747.1429000000001
$ws.Range("K74").Value = 747.1429000000001
$ws.Range("M74").Value = 126.8570999999999
$ws.Range("H77").Value = 925.9655
$ws.Range("I77").Value = 747.1429000000001
$ws.Range("K77").Value = 3735.7145
$ws.Range("M77").Value = 632.2855
$ws.Range("H116").Value = 2602.2222
$ws.Range("I116").Value = 2593.7646
$ws.Range("K116").Value = 2593.7646
$ws.Range("M116").Value = -299.7646
$ws.Range("H122").Value = 7917.237
$ws.Range("I122").Value = 7921.6294
$ws.Range("J122").Value = 7906.4546
$ws.Range("K122").Value = 23764.8882
$ws.Range("L122").Value = 23719.3638
$ws.Range("M122").Value = -21314.8882
$ws.Range("N122").Value = -28619.3638
$ws.Range("H132").Value = 3392.5435
$ws.Range("I132").Value = 1519.6052
$ws.Range("J132").Value = 12289
$ws.Range("K132").Value = 4558.8156
$ws.Range("L132").Value = 36867
$ws.Range("M132").Value = -2028.8156
$ws.Range("N132").Value = -41927

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2602.2222
$ws.Range("I3").Value = 2593.7646
$ws.Range("K3").Value = 2593.7646
$ws.Range("M3").Value = -2479.7646
$ws.Range("H20").Value = 1687.5
$ws.Range("I20").Value = 1622.9
$ws.Range("K20").Value = 1622.9
$ws.Range("M20").Value = -1375.9
$ws.Range("H80").Value = 806.53845
$ws.Range("J80").Value = 967.75
$ws.Range("L80").Value = 967.75
$ws.Range("N80").Value = -2963.75
$ws.Range("H82").Value = 9121.817999999999
$ws.Range("J82").Value = 56708
$ws.Range("L82").Value = 56708
$ws.Range("N82").Value = -57474
$ws.Range("H83").Value = 806.53845
$ws.Range("J83").Value = 967.75
$ws.Range("L83").Value = 4838.75
$ws.Range("N83").Value = -14822.75
$ws.Range("H85").Value = 9121.817999999999
$ws.Range("J85").Value = 56708
$ws.Range("L85").Value = 56708
$ws.Range("N85").Value = -59360
$ws.Range("H99").Value = 1880
$ws.Range("I99").Value = 1880
$ws.Range("K99").Value = 1880
$ws.Range("M99").Value = -382
$ws.Range("H134").Value = 2146.0952
$ws.Range("I134").Value = 1544.75
$ws.Range("K134").Value = 4634.25
$ws.Range("M134").Value = -2099.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1706.0769
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 7500
$ws.Range("N132").Value = -12560
$ws.Range("H134").Value = 1846.9667
$ws.Range("I134").Value = 1899.6364
$ws.Range("K134").Value = 5698.9092
$ws.Range("M134").Value = -3163.9092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9395170
$ws.Range("I4").Value = 11666983
$ws.Range("K4").Value = 35000949
$ws.Range("M4").Value = -35000837
$ws.Range("H34").Value = 2006.1428
$ws.Range("I34").Value = 2006.1428
$ws.Range("K34").Value = 6018.428400000001
$ws.Range("M34").Value = -5934.428400000001
$ws.Range("H44").Value = 67455.734
$ws.Range("I44").Value = 867.25
$ws.Range("K44").Value = 2601.75
$ws.Range("M44").Value = -2203.75
$ws.Range("H93").Value = 500449.5
$ws.Range("I93").Value = 899
$ws.Range("J93").Value = 1000000
$ws.Range("K93").Value = 2697
$ws.Range("L93").Value = 3000000
$ws.Range("M93").Value = -825
$ws.Range("N93").Value = -3003744
$ws.Range("H95").Value = 11920
$ws.Range("J95").Value = 11920
$ws.Range("L95").Value = 35760
$ws.Range("N95").Value = -39878
$ws.Range("H124").Value = 4913.1665
$ws.Range("J124").Value = 11676.5
$ws.Range("L124").Value = 35029.5
$ws.Range("N124").Value = -44849.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 18102.666
$ws.Range("J49").Value = 18102.666
$ws.Range("L49").Value = 18102.666
$ws.Range("N49").Value = -18470.666
$ws.Range("H70").Value = 8017.8335
$ws.Range("I70").Value = 3775.875
$ws.Range("K70").Value = 3775.875
$ws.Range("M70").Value = -3505.875
$ws.Range("H73").Value = 8017.8335
$ws.Range("I73").Value = 3775.875
$ws.Range("K73").Value = 3775.875
$ws.Range("M73").Value = -2839.875
$ws.Range("H103").Value = 85034
$ws.Range("J103").Value = 85034
$ws.Range("L103").Value = 85034
$ws.Range("N103").Value = -87378
$ws.Range("H122").Value = 15004
$ws.Range("I122").Value = 15004
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 45012
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -42562
$ws.Range("N122").ClearContents()
$ws.Range("H124").Value = 78221.44500000001
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H132").Value = 7257.3
$ws.Range("I132").Value = 7431.489
$ws.Range("J132").Value = 5689.6
$ws.Range("K132").Value = 22294.467
$ws.Range("L132").Value = 17068.8
$ws.Range("M132").Value = -19764.467
$ws.Range("N132").Value = -22128.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4846.973
$ws.Range("I40").Value = 4356.033
$ws.Range("K40").Value = 4356.033
$ws.Range("M40").Value = -4220.033
$ws.Range("H41").Value = 29161.8
$ws.Range("I41").Value = 29161.8
$ws.Range("K41").Value = 29161.8
$ws.Range("M41").Value = -28723.8
$ws.Range("H55").Value = 1513.7059
$ws.Range("I55").Value = 268
$ws.Range("J55").Value = 2915.125
$ws.Range("K55").Value = 268
$ws.Range("L55").Value = 2915.125
$ws.Range("M55").Value = -95
$ws.Range("N55").Value = -3261.125
$ws.Range("H82").Value = 21450.143
$ws.Range("I82").Value = 18358.5
$ws.Range("J82").Value = 40000
$ws.Range("K82").Value = 18358.5
$ws.Range("L82").Value = 40000
$ws.Range("M82").Value = -17997.5
$ws.Range("N82").Value = -40722
$ws.Range("H85").Value = 21450.143
$ws.Range("I85").Value = 18358.5
$ws.Range("J85").Value = 40000
$ws.Range("K85").Value = 18358.5
$ws.Range("L85").Value = 40000
$ws.Range("M85").Value = -17110.5
$ws.Range("N85").Value = -42496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 159426.4
$ws.Range("I62").Value = 262377.34
$ws.Range("K62").Value = 262377.34
$ws.Range("M62").Value = -261753.34
$ws.Range("H65").Value = 159426.4
$ws.Range("I65").Value = 262377.34
$ws.Range("K65").Value = 1311886.7
$ws.Range("M65").Value = -1308766.7
$ws.Range("H106").Value = 44870.5
$ws.Range("J106").Value = 44870.5
$ws.Range("L106").Value = 44870.5
$ws.Range("N106").Value = -47394.5
$ws.Range("H124").Value = 200214
$ws.Range("J124").Value = 200214
$ws.Range("L124").Value = 200214
$ws.Range("N124").Value = -210034
